$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the front; existing columns (and their shared-string
# values / formatting) shift from A:D to B:E automatically.
$ws.Columns.Item(1).Insert()

# New column A header/value
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# Updated Neo4j queries (B2 = cases query, C2 = stat query) - these cells
# held the old query text before the column shift; overwrite with new text.
$ws.Range("B2").Value = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)`n    WHERE c.race = `"BLACK_OR_AFRICAN_AMERICAN`"`nWITH DISTINCT c, a, ct`nRETURN `n    COALESCE(c.case_id, '') AS ``Case ID``,`n    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n    COALESCE(a.arm_id, '') AS ``Arm``,`n    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,`n    COALESCE(c.disease, '') AS ``Diagnosis``,`n    COALESCE(c.gender, '') AS ``Gender``,`n    COALESCE(c.race, '') AS ``Race``,`n    COALESCE(c.ethnicity, '') AS ``Ethnicity``"

$ws.Range("C2").Value = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n    WHERE c.race = `"BLACK_OR_AFRICAN_AMERICAN`"`nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`nRETURN `n    COUNT(DISTINCT f) AS number_of_files,`n    COUNT(DISTINCT c.case_id) AS number_of_cases,`n    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

# Formatting: wrap text on B2/C2 (already inherited from old A2/B2 style),
# row 2 taller to fit the now-longer wrapped text, and a narrow best-fit
# column A.
$ws.Rows.Item(2).RowHeight = 174
$ws.Columns.Item(1).ColumnWidth = 7.92
$ws.Range("B2:C2").WrapText = $true

# Restore the selection Excel leaves after this kind of edit.
$ws.Range("B5").Select() | Out-Null
